# Auto-generated edit script: updates scheduled market price/profit data
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 2478.4783
$ws.Range("I39").Value = 411.35715
$ws.Range("J39").Value = 5694
$ws.Range("K39").Value = 1234.07145
$ws.Range("L39").Value = 17082
$ws.Range("M39").Value = -938.0714499999999
$ws.Range("N39").Value = -17674
$ws.Range("H45").Value = 3361
$ws.Range("I45").Value = 3375
$ws.Range("J45").Value = 3333
$ws.Range("K45").Value = 10125
$ws.Range("L45").Value = 9999
$ws.Range("M45").Value = -9933
$ws.Range("N45").Value = -10383
$ws.Range("H53").Value = 164.71428
$ws.Range("I53").Value = 92.2
$ws.Range("K53").Value = 92.2
$ws.Range("M53").Value = 544.8
$ws.Range("H58").Value = 676
$ws.Range("I58").Value = 460
$ws.Range("K58").Value = 1380
$ws.Range("M58").Value = -1230
$ws.Range("H95").Value = 36824.6
$ws.Range("J95").Value = 36824.6
$ws.Range("L95").Value = 36824.6
$ws.Range("N95").Value = -42316.6
$ws.Range("H112").Value = 5396.909
$ws.Range("I112").Value = 1500
$ws.Range("K112").Value = 4500
$ws.Range("M112").Value = -3392
$ws.Range("H138").Value = 5871.4
$ws.Range("J138").Value = 7090.9546
$ws.Range("L138").Value = 21272.8638
$ws.Range("N138").Value = -31552.8638
$ws.Range("H141").Value = 599.8889
$ws.Range("I141").Value = 599.8889
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 1799.6667
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 3380.3333
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2746.1558
$ws.Range("I32").Value = 572.25714
$ws.Range("K32").Value = 572.25714
$ws.Range("M32").Value = -285.25714
$ws.Range("H132").Value = 2569.8809
$ws.Range("I132").Value = 2152.3076
$ws.Range("K132").Value = 6456.9228
$ws.Range("M132").Value = -3926.9228

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4056.8845
$ws.Range("I20").Value = 4182.6875
$ws.Range("J20").Value = 3855.6
$ws.Range("K20").Value = 4182.6875
$ws.Range("L20").Value = 3855.6
$ws.Range("M20").Value = -3935.6875
$ws.Range("N20").Value = -4349.6
$ws.Range("H86").Value = 6411799.5
$ws.Range("I86").Value = 9260499
$ws.Range("J86").Value = 2224.75
$ws.Range("K86").Value = 9260499
$ws.Range("L86").Value = 2224.75
$ws.Range("M86").Value = -9259376
$ws.Range("N86").Value = -4470.75
$ws.Range("H89").Value = 6411799.5
$ws.Range("I89").Value = 9260499
$ws.Range("J89").Value = 2224.75
$ws.Range("K89").Value = 46302495
$ws.Range("L89").Value = 11123.75
$ws.Range("M89").Value = -46296879
$ws.Range("N89").Value = -22355.75
$ws.Range("H94").Value = 2520.2778
$ws.Range("I94").Value = 1418.6316
$ws.Range("J94").Value = 3751.5293
$ws.Range("K94").Value = 1418.6316
$ws.Range("L94").Value = 3751.5293
$ws.Range("M94").Value = -967.6315999999999
$ws.Range("N94").Value = -4653.5293
$ws.Range("H99").Value = 3095.7368
$ws.Range("I99").Value = 2966.2144
$ws.Range("K99").Value = 2966.2144
$ws.Range("M99").Value = -1468.2144
$ws.Range("H134").Value = 3848.4243
$ws.Range("I134").Value = 3486.3794
$ws.Range("J134").Value = 6473.25
$ws.Range("K134").Value = 10459.1382
$ws.Range("L134").Value = 19419.75
$ws.Range("M134").Value = -7924.138199999999
$ws.Range("N134").Value = -24489.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 2500
$ws.Range("I8").Value = 2500
$ws.Range("K8").Value = 2500
$ws.Range("M8").Value = -2360
$ws.Range("H31").Value = 5338.793
$ws.Range("I31").Value = 5589.6587
$ws.Range("J31").Value = 4733.7646
$ws.Range("K31").Value = 5589.6587
$ws.Range("L31").Value = 4733.7646
$ws.Range("M31").Value = -5294.6587
$ws.Range("N31").Value = -5323.7646
$ws.Range("H34").Value = 5338.793
$ws.Range("I34").Value = 5589.6587
$ws.Range("J34").Value = 4733.7646
$ws.Range("K34").Value = 5589.6587
$ws.Range("L34").Value = 4733.7646
$ws.Range("M34").Value = -5387.6587
$ws.Range("N34").Value = -5137.7646
$ws.Range("H122").Value = 216608.14
$ws.Range("I122").Value = 302868.4
$ws.Range("J122").Value = 957.5
$ws.Range("K122").Value = 908605.2000000001
$ws.Range("L122").Value = 2872.5
$ws.Range("M122").Value = -906155.2000000001
$ws.Range("N122").Value = -7772.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1589.5385
$ws.Range("I107").Value = 1831.1428
$ws.Range("J107").Value = 1307.6666
$ws.Range("K107").Value = 5493.428400000001
$ws.Range("L107").Value = 3922.9998
$ws.Range("M107").Value = -3573.428400000001
$ws.Range("N107").Value = -7762.9998
$ws.Range("H131").Value = 2006.1154
$ws.Range("I131").Value = 1409.0834
$ws.Range("K131").Value = 4227.2502
$ws.Range("M131").Value = 812.7497999999996
$ws.Range("H140").Value = 912.8684
$ws.Range("I140").Value = 912.8684
$ws.Range("K140").Value = 2738.6052
$ws.Range("M140").Value = 2441.3948

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6848.25
$ws.Range("I70").Value = 6233.3335
$ws.Range("J70").Value = 7053.222
$ws.Range("K70").Value = 6233.3335
$ws.Range("L70").Value = 7053.222
$ws.Range("M70").Value = -5963.3335
$ws.Range("N70").Value = -7593.222
$ws.Range("H73").Value = 6848.25
$ws.Range("I73").Value = 6233.3335
$ws.Range("J73").Value = 7053.222
$ws.Range("K73").Value = 6233.3335
$ws.Range("L73").Value = 7053.222
$ws.Range("M73").Value = -5297.3335
$ws.Range("N73").Value = -8925.222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3388.5
$ws.Range("I22").Value = 3388.5
$ws.Range("K22").Value = 3388.5
$ws.Range("M22").Value = -3093.5
$ws.Range("H27").Value = 3388.5
$ws.Range("I27").Value = 3388.5
$ws.Range("K27").Value = 3388.5
$ws.Range("M27").Value = -3281.5
$ws.Range("H132").Value = 3008.4348
$ws.Range("I132").Value = 2344.2222
$ws.Range("K132").Value = 7032.6666
$ws.Range("M132").Value = -4502.6666
$ws.Range("H136").Value = 2579.45
$ws.Range("I136").Value = 2615.2104
$ws.Range("K136").Value = 7845.6312
$ws.Range("M136").Value = -5295.6312

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H101").Value = 40665.668
$ws.Range("J101").Value = 29998.5
$ws.Range("L101").Value = 29998.5
$ws.Range("N101").Value = -36488.5
$ws.Range("H113").Value = 1205.3
$ws.Range("I113").Value = 1299.4
$ws.Range("J113").Value = 1111.2
$ws.Range("K113").Value = 3898.2
$ws.Range("L113").Value = 3333.6
$ws.Range("M113").Value = -1728.2
$ws.Range("N113").Value = -7673.6
$ws.Range("H132").Value = 3383.8
$ws.Range("I132").Value = 3383.8
$ws.Range("K132").Value = 10151.4
$ws.Range("M132").Value = -7621.400000000001

Write-Host "Applied 177 cell updates."
